# Realestate Update resale numbers 2024-01-15 22:20
# Appends a new data row (row 61) to the CityResaleNum sheet with the
# latest resale number snapshot, mirroring the structure of existing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 61

# Text columns (Date, Time, Weekday, Week) - force as text so values like
# "02" and the date/time strings are not reinterpreted as numbers/dates.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 2).NumberFormat = "@"
$ws.Cells.Item($row, 3).NumberFormat = "@"
$ws.Cells.Item($row, 4).NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2024-01-15"
$ws.Cells.Item($row, 2).Value = "22:20:51"
$ws.Cells.Item($row, 3).Value = "Monday"
$ws.Cells.Item($row, 4).Value = "02"

# Numeric columns (city resale numbers)
$ws.Cells.Item($row, 5).Value = 138862
$ws.Cells.Item($row, 6).Value = 139111
$ws.Cells.Item($row, 7).Value = 171463
$ws.Cells.Item($row, 8).Value = 148097
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 119408
$ws.Cells.Item($row, 11).Value = 220996
$ws.Cells.Item($row, 12).Value = 254086
$ws.Cells.Item($row, 13).Value = 184945
$ws.Cells.Item($row, 14).Value = 110403
$ws.Cells.Item($row, 15).Value = 41150
$ws.Cells.Item($row, 16).Value = 30891
$ws.Cells.Item($row, 17).Value = 73194
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 42639
$ws.Cells.Item($row, 20).Value = -1
